$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '28.533.18'
Set-TextValue $ws.Range('E2') '  +1.76%  '
Set-TextValue $ws.Range('D3') '1.864.71'
Set-TextValue $ws.Range('E3') '  +1.84%  '
Set-TextValue $ws.Range('E4') '  +0.18%  '
Set-TextValue $ws.Range('D5') '323.71'
Set-TextValue $ws.Range('E5') '  -0.41%  '
Set-TextValue $ws.Range('D6') '1.002'
Set-TextValue $ws.Range('E6') '  +0.16%  '
Set-TextValue $ws.Range('D7') '0.4596'
Set-TextValue $ws.Range('E7') '  -1.32%  '
Set-TextValue $ws.Range('E8') '  +0.03%  '
Set-TextValue $ws.Range('D10') '0.9722'
Set-TextValue $ws.Range('E10') '  +1.32%  '
Set-TextValue $ws.Range('D11') '21.88'
Set-TextValue $ws.Range('E11') '  +0.01%  '
Set-TextValue $ws.Range('D12') '1.862.05'
Set-TextValue $ws.Range('E12') '  +1.73%  '
Set-TextValue $ws.Range('D13') '6.969'
Set-TextValue $ws.Range('E13') '  +0.93%  '
Set-TextValue $ws.Range('D14') '5.674'
Set-TextValue $ws.Range('E14') '  -0.03%  '
Set-TextValue $ws.Range('D15') '0.06935'
Set-TextValue $ws.Range('E15') '  +1.09%  '
Set-TextValue $ws.Range('D16') '88.05'
Set-TextValue $ws.Range('E16') '  +0.96%  '
Set-TextValue $ws.Range('D17') '1.004'
Set-TextValue $ws.Range('E17') '  +0.30%  '
Set-TextValue $ws.Range('D18') '0.000009984'
Set-TextValue $ws.Range('E18') '  +0.68%  '
Set-TextValue $ws.Range('E19') '  +1.10%  '
Set-TextValue $ws.Range('E20') '  +0.12%  '
Set-TextValue $ws.Range('D21') '28.537.71'
Set-TextValue $ws.Range('E21') '  +1.78%  '
Set-TextValue $ws.Range('D22') '5.250'
Set-TextValue $ws.Range('E22') '  -1.27%  '
Set-TextValue $ws.Range('D23') '11.05'
Set-TextValue $ws.Range('E23') '  +0.67%  '
Set-TextValue $ws.Range('D24') '2.105'
Set-TextValue $ws.Range('E24') '  +0.72%  '
Set-TextValue $ws.Range('D25') '2.080.58'
Set-TextValue $ws.Range('E25') '  +0.31%  '
Set-TextValue $ws.Range('D26') '152.28'
Set-TextValue $ws.Range('E26') '  -1.01%  '
Set-TextValue $ws.Range('D27') '19.23'
Set-TextValue $ws.Range('E27') '  +0.68%  '
Set-TextValue $ws.Range('D28') '5.773'
Set-TextValue $ws.Range('E28') '  +1.43%  '
Set-TextValue $ws.Range('E29') '  +1.12%  '
Set-TextValue $ws.Range('D30') '119.05'
Set-TextValue $ws.Range('E30') '  +1.15%  '
Set-TextValue $ws.Range('D31') '0.09303'
Set-TextValue $ws.Range('E31') '  +0.51%  '
Set-TextValue $ws.Range('D32') '0.9109'
Set-TextValue $ws.Range('E32') '  -2.67%  '
Set-TextValue $ws.Range('D33') '5.258'
Set-TextValue $ws.Range('E33') '  -0.37%  '
Set-TextValue $ws.Range('D34') '1.326'
Set-TextValue $ws.Range('E34') '  +0.56%  '
Set-TextValue $ws.Range('D35') '3.319'
Set-TextValue $ws.Range('E35') '  +0.86%  '
Set-TextValue $ws.Range('D36') '0.05772'
Set-TextValue $ws.Range('E36') '  -1.48%  '
Set-TextValue $ws.Range('D37') '1.146'
Set-TextValue $ws.Range('E37') '  +0.45%  '
Set-TextValue $ws.Range('D38') '0.02073'
Set-TextValue $ws.Range('E38') '  -2.56%  '
Set-TextValue $ws.Range('D39') '7.667'
Set-TextValue $ws.Range('E39') '  -1.87%  '
Set-TextValue $ws.Range('D40') '0.5606'
Set-TextValue $ws.Range('E40') '  +0.44%  '
Set-TextValue $ws.Range('D41') '0.1776'
Set-TextValue $ws.Range('E41') '  +1.06%  '
Set-TextValue $ws.Range('D42') '9.730'
Set-TextValue $ws.Range('E42') '  -1.43%  '
Set-TextValue $ws.Range('D43') '0.07201'
Set-TextValue $ws.Range('E43') '  +2.70%  '
Set-TextValue $ws.Range('D44') '11.65'
Set-TextValue $ws.Range('E44') '  +0.58%  '
Set-TextValue $ws.Range('D45') '0.5267'
Set-TextValue $ws.Range('E45') '  +0.10%  '
Set-TextValue $ws.Range('D46') '2.145'
Set-TextValue $ws.Range('E46') '  +0.77%  '
Set-TextValue $ws.Range('D47') '1.137'
Set-TextValue $ws.Range('E47') '  +1.75%  '
Set-TextValue $ws.Range('E48') '  +0.08%  '
Set-TextValue $ws.Range('D49') '112.35'
Set-TextValue $ws.Range('E49') '  -0.39%  '
Set-TextValue $ws.Range('D50') '2.414'
Set-TextValue $ws.Range('E50') '  +4.16%  '
Set-TextValue $ws.Range('E51') '  +0.19%  '
